# Fruta / hortaliza, semanal
# Re-shuffles the per-record columns (Fecha, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Precio $/Kg) across the
# existing data rows (2..41) of the sheet. Every other column (Mercado,
# Region, Codreg, Tipo, Producto, Categoria, Variedad, Calidad, Unidad de
# comercializacion, Origen, Kg/unidad) stays put; only the row that a
# given (date, volume, prices) tuple lives on changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# targetRow -> sourceRow : the data that ends up on $targetRow is the data
# that used to live on $sourceRow (before this edit was applied).
$rowMap = @{
    2  = 6
    3  = 10
    4  = 2
    5  = 32
    6  = 36
    7  = 14
    8  = 23
    9  = 34
    10 = 5
    11 = 28
    12 = 12
    13 = 35
    14 = 19
    15 = 8
    16 = 29
    17 = 38
    18 = 40
    19 = 30
    20 = 26
    21 = 15
    22 = 21
    23 = 3
    24 = 17
    25 = 16
    26 = 22
    27 = 7
    28 = 24
    29 = 27
    30 = 20
    31 = 41
    32 = 39
    33 = 11
    34 = 4
    35 = 9
    36 = 37
    37 = 18
    38 = 31
    39 = 33
    40 = 25
    41 = 13
}

$cols = @("D", "M", "N", "O", "P", "S")

# Snapshot every affected cell's current value before writing anything,
# since the remap is an in-place permutation (several rows read from each
# other).
$snapshot = @{}
foreach ($row in $rowMap.Keys) {
    $srcRow = $rowMap[$row]
    foreach ($col in $cols) {
        $key = "$col$srcRow"
        if (-not $snapshot.ContainsKey($key)) {
            $snapshot[$key] = $ws.Range($key).Value()
        }
    }
}

foreach ($row in $rowMap.Keys) {
    $srcRow = $rowMap[$row]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $snapshot["$col$srcRow"]
    }
}
